$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark. It currently sits on the
#    last (empty) paragraph of the document; after this edit it will
#    live on the newly inserted paragraph instead.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Insert a brand-new list paragraph right after the
#    "Didn't include tests of the user input..." bullet. Using
#    InsertParagraphAfter() on that paragraph's Range makes the new
#    paragraph inherit the same pPr (ListParagraph style, numPr for
#    numId 1, spacing) automatically.
# ------------------------------------------------------------------
$sourcePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*include tests of the user input based methods*") {
        $sourcePara = $p
    }
}

$sourcePara.Range.InsertParagraphAfter()
$newPara = $sourcePara.Next()

$r = $newPara.Range
$r.Collapse(1)

# ------------------------------------------------------------------
# 3) Type the paragraph text as two separate runs (matching the
#    target markup) by placing a temporary bookmark between them --
#    that forces the engine to keep them as distinct <w:r> elements
#    instead of silently merging identically formatted runs.
# ------------------------------------------------------------------
$r.InsertAfter("Didn" + [char]0x2019 + "t include tests of the ")
$r.Collapse(0)
$r.Bookmarks.Add("zzzTempRunSplit") | Out-Null
$r.InsertAfter("file IO " + [char]0x2013 + " we just hardcoded a pack in the test code.")
$d.Bookmarks.Item("zzzTempRunSplit").Delete()

# ------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark at the very end of the new
#    paragraph's text (right after the two runs, before the
#    paragraph mark). A collapsed bookmark placed directly at that
#    boundary is mishandled by the runtime, so instead: append a
#    throw-away character, wrap a (non-collapsed) bookmark around
#    it, then delete the character -- the bookmark gracefully
#    collapses to the correct position.
# ------------------------------------------------------------------
$r.Collapse(0)
$r.InsertAfter("x")
$d.Bookmarks.Add("_GoBack", $r)
$r.Text = ""
